$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 76; this shifts the existing rows 76-96 down to 77-97
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new record's data
$ws.Cells.Item(76, 1).Value = 10
$ws.Cells.Item(76, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(76, 3).Value = "La Araucanía"
$ws.Cells.Item(76, 4).Value = 44855
$ws.Cells.Item(76, 5).Value = 9
$ws.Cells.Item(76, 6).Value = 100112022
$ws.Cells.Item(76, 7).Value = "Arveja Verde"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 30
$ws.Cells.Item(76, 11).Value = 19000
$ws.Cells.Item(76, 12).Value = 20000
$ws.Cells.Item(76, 13).Value = 19333
$ws.Cells.Item(76, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(76, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(76, 16).Value = 773
$ws.Cells.Item(76, 17).Value = 25
$ws.Cells.Item(76, 18).Value = "Hortaliza"
